$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("experiments")

# Row 2
$ws.Range("F2").Value = 52
$ws.Range("G2").Value = 2117
$ws.Range("H2").Value = 9
$ws.Range("I2").Value = 2186.201
$ws.Range("J2").Value = 48462.919
$ws.Range("K2").Value = 215.409

# Row 3
$ws.Range("F3").Value = 54
$ws.Range("G3").Value = 2218
$ws.Range("H3").Value = 8
$ws.Range("I3").Value = 2289.898
$ws.Range("J3").Value = 52099.638
$ws.Range("K3").Value = 207.739

# Row 4
$ws.Range("F4").Value = 50
$ws.Range("G4").Value = 21342
$ws.Range("H4").Value = 20
$ws.Range("I4").Value = 21421.827
$ws.Range("J4").Value = 535151.941
$ws.Range("K4").Value = 670.454
